$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns are treated as literal text,
# matching the original inline-string cell type, so values such as
# "4.20" or "0.0273" keep their exact digits/trailing zeros instead of
# being auto-coerced to numbers by the COM Value setter.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.744.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.174.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.87%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.66%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.173.87"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.80%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.727.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.72%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.25"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.670.11"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.177.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.70%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.26%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.310.60"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.83%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.08%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.08%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.11"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.79%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.43"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.46"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.43%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.646.66"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "327.46"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0273"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.80%  "
